$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = 100,101,102,103,104,105,106,107,108,109,110,111
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("C$row")
    $cell.Value = $values[$i]
    $cell.Font.Color = 255
}

$ws.Range("B17").Select()
